# Applies the commit "Added enums + drivers + a lot of refactoring":
#  1. Rename sheet "LogConfigs" -> "LoginConfigs"
#  2. Update DriverConfigs driver-location values (geckodriver / msedgedriver)
#  3. Widen column E on DriverConfigs sheet
#  4. Restyle E2/F2 on DriverConfigs sheet (drop explicit Arial font -> default font)

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item(1)
$wsDriver = $wb.Worksheets.Item(2)

# 1. Rename the first sheet.
$wsLogin.Name = "LoginConfigs"

# 2. Update the driver executable location values.
$wsDriver.Range("E2").Value = "src/main/resources/geckodriver.exe"
$wsDriver.Range("F2").Value = "src/main/resources/msedgedriver.exe"

# 4. Match the style used elsewhere (e.g. LoginConfigs!E1) for E2:F2, which
#    uses the plain/default font instead of the explicit Arial font. Using
#    copy/PasteSpecial(formats) keeps the cell values untouched while only
#    bringing over the formatting (font/fill/etc.).
$wsLogin.Range("E1").Copy()
$wsDriver.Range("E2:F2").PasteSpecial(-4122)

# 3. Widen column E. The COM layer quantizes ColumnWidth to pixel units, so
#    31.1 is the input that lands exactly on the 32-wide pixel boundary.
$wsDriver.Columns.Item(5).ColumnWidth = 31.1

$excel.CutCopyMode = $false
